$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format first,
# matching the original inlineStr (text) cell type in the workbook.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.373.69'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.874.48'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '0.7143'
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.3112'
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("D9").Value = '0.07771'
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("D10").Value = '25.15'
$ws.Range("E10").Value = '  +1.76%  '
$ws.Range("D11").Value = '0.08448'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = '1.869.02'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '5.254'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '0.7136'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '91.24'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '29.375.32'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '6.099'
$ws.Range("E17").Value = '  +2.99%  '
$ws.Range("D18").Value = '0.000008243'
$ws.Range("E18").Value = '  +5.39%  '
$ws.Range("D19").Value = '240.96'
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("D21").Value = '2.124.96'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '7.779'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '0.1596'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '163.29'
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D27").Value = '9.069'
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("D28").Value = '18.55'
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").Value = '1.514'
$ws.Range("E29").Value = '  +1.02%  '
$ws.Range("D30").Value = '4.422'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.335'
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '1.292'
$ws.Range("E32").Value = '  -1.83%  '
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").Value = '0.7439'
$ws.Range("E36").Value = '  -7.49%  '
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").Value = '0.01872'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").Value = '1.226.93'
$ws.Range("E39").Value = '  +4.54%  '
$ws.Range("D40").Value = '2.729'
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("D41").Value = '6.525'
$ws.Range("E41").Value = '  +5.65%  '
$ws.Range("D42").Value = '110.83'
$ws.Range("E42").Value = '  +8.97%  '
$ws.Range("D43").Value = '0.8889'
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("D44").Value = '72.89'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = '2.021.93'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D47").Value = '1.815'
$ws.Range("E47").Value = '  +2.57%  '
$ws.Range("D48").Value = '0.5214'
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").Value = '0.00000000122'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").Value = '9.447'
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").Value = '0.4323'
$ws.Range("E51").Value = '  +1.52%  '
